$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before H; this shifts existing H..X to I..Y
$ws.Columns("H").Insert()

# The new column should have the same (~28.44 char) width as its neighbours
$ws.Columns("H").ColumnWidth = 27.65

# Populate the new header cell (row 1) - reuses the default header style
$ws.Range("H1").Value = "personal_account_number2"

# Populate the new data cell (row 2) with the new account number, matching the
# wrapped-text style used by the neighbouring account-number/name cells
$ws.Range("H2").Value = "RS35 2059 0310 0441 7882 84"
$ws.Range("H2").WrapText = $true

# Row 3 has no data in the new column - remove the cell entirely
$ws.Range("H3").Clear()

# Update the view: scroll position and selection
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("H1").Select()
